$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("AK2").Value = 9

# Row 3 changes
$ws.Range("G3").Value = 1.4
$ws.Range("H3").Value = 4.5
$ws.Range("I3").Value = 7.5
$ws.Range("L3").Value = 7
$ws.Range("N3").Value = 12
$ws.Range("W3").Value = 7
$ws.Range("AD3").Value = 8.5
$ws.Range("AK3").Value = 81
$ws.Range("AO3").Value = 7
$ws.Range("AT3").Value = 9
$ws.Range("AV3").Value = 8.5
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 301
